$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (engine quantizes ColumnWidth to 1/6-character steps,
# so we target the bucket midpoint closest to the desired stored width)
$ws.Columns.Item(8).ColumnWidth = 2.3333333333333335
$ws.Columns.Item(11).ColumnWidth = 4.833333333333333

# Update row 1 values
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 19
$ws.Range("E1").Value = 3
$ws.Range("F1").Value = 2
$ws.Range("G1").Value = 17
$ws.Range("H1").Value = 10
$ws.Range("I1").Value = 23
$ws.Range("J1").Value = 24
$ws.Range("K1").Value = 0.064
$ws.Range("L1").Value = 0.064
$ws.Range("M1").Value = 0.036
$ws.Range("N1").Value = 0.037
